# Applies the PressRelease.docx text edits described in the commit diff.
# Uses Find/Replace across the Word object model so Word's own run
# splitting/merging produces the final OOXML runs.

$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $range = $d.Content
    $ok = $range.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2)
    if (-not $ok) {
        Write-Output "NOT FOUND: $find"
    }
}

# --- Paragraph: "El suscriptor ..." ---
# Merge "carga " + "layouts" (spell-checked) + " y " into one run, dropping
# the proofErr spell-check markers around "layouts".
Replace-Text "carga layouts y COCO" "carga layouts y COCO"

# "Ofrece los siguientes servicios:" -> "COCO ofrece los siguientes servicios:"
Replace-Text "procesa, almacena y transforma los datos. Ofrece los siguientes servicios:" "procesa, almacena y transforma los datos. COCO ofrece los siguientes servicios:"

# --- Paragraph: "Cálculo de la tarifa de línea" ---
Replace-Text "contiene un algoritmo capaz de calcular la prima de pólizas experiencia global con base en las condiciones de la póliza seleccionadas. Además, almacena el cálculo realizado" "COCO contiene un algoritmo capaz de calcular la prima de pólizas “experiencia global” con base en las condiciones de la póliza seleccionadas. Además, Coco almacena el cálculo realizado"

# --- Paragraph: "Generación del formato de cotización" ---
Replace-Text ": genera un formato de cotización con las condiciones de la póliza seleccionadas y la prima calculada" ": Coco genera un formato de cotización con las condiciones de la póliza seleccionadas y la prima calculada"

# --- Paragraph: Santiago Redondo quote (first half) ---
Replace-Text "“A pesar de cumplir con las metas de crecimiento de cartera en los últimos 3 años, no hemos visto un aumento en el crecimiento de la cartera de pólizas experiencia global. Es fundamental aumentar el primaje de estas cuentas, pues no podemos depender tanto de las pólizas grandes negocios para llegar a la meta. Estas pueden cambiar de compañía en cualquier momento, incluso brindándoles un excelente servicio. Su política interna les exige cambiar cada cierto tiempo de aseguradora. " "“A pesar de cumplir con las metas de crecimiento de cartera en los últimos 3 años, no hemos visto un aumento en el crecimiento de la cartera de pólizas “experiencia global”. Es fundamental aumentar el primaje de estas cuentas, pues no podemos depender tanto de las pólizas de grandes negocios para llegar a la meta ya que estas pueden cambiar de compañía en cualquier momento, incluso brindándoles un excelente servicio y mejor que el de los competidores. En algunos casos, por ejemplo, existen políticas internas que les exigen a las compañías cambiar cada cierto tiempo de aseguradora. "

# --- Paragraph: Santiago Redondo quote (second half) ---
Replace-Text " automatizar el proceso de suscripción de pólizas experiencia global vamos a disminuir el tiempo de respuesta al cliente, liberar espacio para poder usar de mejor forma el sistema central de la compañía y recopilar información para identificar las regiones en donde debemos impulsar la captación de nuevos agentes. Además, los suscriptores tendrán más tiempo para realizar un análisis exhaustivo a las nuevas cuentas experiencia global”. " " automatizar el proceso de suscripción de pólizas “experiencia global”, vamos a disminuir el tiempo de respuesta al cliente, liberar espacio para poder usar de mejor forma el sistema central de la compañía y recopilar información para identificar las regiones en donde debemos impulsar la captación de nuevos agentes. Además, los suscriptores tendrán más tiempo para realizar un análisis exhaustivo a las nuevas cuentas “experiencia global”. "

# --- Paragraph: "Ha sido fundamental ..." ---
Replace-Text "Ha sido fundamental el trabajo de distintas áreas, suscripción, siniestros, ventas y emisión, para que " "Ha sido fundamental el trabajo de distintas áreas como Suscripción, Siniestros, Ventas y Emisión, para que "

Write-Output "done"
